# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion summary text in A1 ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.32 = 12721.13 pesos`n✅ 12721.13 pesos = 3.31 = 971.31 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$hoja1.Range("A1").Value = $newText

# --- tasas: update the N10/O10/N12/O12 rate figures ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 301
$tasas.Range("O10").Value = 3829.06
$tasas.Range("N12").Value = 3840
$tasas.Range("O12").Value = 293.2
